# "updated main GSC export data"
#
# The Google Search Console "Video Indexing" export gained four more days
# of data (2025-11-28 .. 2025-12-01) on the "Chart" sheet, and the rollup
# "Failed" video count on the "Table" sheet dropped from 24 to 23 now that
# one more video's indexing issue cleared up (the new 2025-12-01 row has no
# "Impressions" value yet, same as how freshly-added days look elsewhere in
# the sheet).

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")
$table = $wb.Worksheets.Item("Table")

# New daily rows appended right after the existing last row (55).
$newRows = @(
    @{ Row = 56; Date = "2025-11-28"; NoVideoIndexed = 23; VideoIndexed = 1; Impressions = 0 },
    @{ Row = 57; Date = "2025-11-29"; NoVideoIndexed = 23; VideoIndexed = 1; Impressions = 0 },
    @{ Row = 58; Date = "2025-11-30"; NoVideoIndexed = 23; VideoIndexed = 1; Impressions = 0 },
    @{ Row = 59; Date = "2025-12-01"; NoVideoIndexed = 23; VideoIndexed = 1; Impressions = $null }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    # Leading "'" forces plain text entry (same as the other Date cells in
    # this column) instead of Excel auto-converting the literal into a date
    # serial number.
    $chart.Range("A$rowNum").Value = "'" + $r.Date
    $chart.Range("B$rowNum").Value = $r.NoVideoIndexed
    $chart.Range("C$rowNum").Value = $r.VideoIndexed
    if ($null -eq $r.Impressions) {
        # Blank/not-yet-available Impressions is stored as an empty string,
        # matching the blank cells used throughout the sheet.
        $chart.Range("D$rowNum").Value = "'"
    } else {
        $chart.Range("D$rowNum").Value = $r.Impressions
    }
}

# Re-apply the sheet's default (unstyled) formatting to the new block so it
# matches the rest of the table instead of keeping any incidental
# text/date formatting picked up while the values were entered.
$chart.Range("A2").Copy()
$chart.Range("A56:D59").PasteSpecial(-4122)

# The "Failed" validation rollup on the Table sheet now reflects one fewer
# failing video (24 -> 23).
$table.Range("C2").Value = 23
